# Final changes for testing part of REF
#
# 1. Duplicate the "Tests" sheet, rename the copy to "Result", and change
#    its header cell B1 from "ExpectedResult" to "ActualResult" (the rest
#    of the data stays identical to "Tests").
# 2. Tighten the data validation list on the new sheet to the used range
#    (B2:B9) instead of the whole column.
# 3. Restore the selection/active-cell state seen in the target workbook:
#    "Tests" ends up with D9 selected (and is no longer the active tab),
#    "Result" ends up the active tab with F12 selected.

$wb  = $excel.ActiveWorkbook
$tests = $wb.Worksheets.Item("Tests")

# Duplicate "Tests" right after itself, then rename + retarget the header.
$tests.Copy($null, $tests) | Out-Null
$result = $wb.Worksheets.Item($tests.Index + 1)
$result.Name = "Result"
$result.Range("B1").Value = "ActualResult"

# Narrow the copied data validation down to the populated rows, keeping
# the same three-option dropdown list.
$result.Range("B2:B1048576").Validation.Delete() | Out-Null
$result.Range("B2:B9").Validation.Add(3, 1, 1, '"Success,BusinessException,SystemException"') | Out-Null

# Match the final selection state: Tests -> D9 selected (inactive tab),
# Result -> F12 selected (active tab, selected last).
$tests.Range("D9").Select() | Out-Null
$result.Range("F12").Select() | Out-Null
